$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column Q (17) = year 2020 data, added to the right of existing P (2019) column.
# Row -> (format source column, value)
# Rows 4,5,13,14 share the same style between P and Q, so copy format from P.
# Rows 6-12 use the "Обычный 2" flavoured style (same as columns D-M) for Q
# instead of the "Обычный" flavoured style used in P, so copy format from D.

$xlPasteFormats = -4122

function Set-Cell($row, $formatSourceCol, $value) {
    $ws.Cells.Item($row, $formatSourceCol).Copy() | Out-Null
    $target = $ws.Cells.Item($row, 17)
    $target.PasteSpecial($xlPasteFormats) | Out-Null
    $target.Value = $value
}

Set-Cell 4 16 2020
Set-Cell 5 16 0.1
Set-Cell 6 4 0.2
Set-Cell 7 4 "-"
Set-Cell 8 4 0.2
Set-Cell 9 4 "-"
Set-Cell 10 4 0.1
Set-Cell 11 4 "-"
Set-Cell 12 4 0.3
Set-Cell 13 16 "-"
Set-Cell 14 16 "-"

$excel.CutCopyMode = $false

# Selection left behind by the editing session, per the commit diff.
$ws.Range("O17").Select()
